# This script rewrites the header row (row 1) of Sheet1 so that it contains
# the new, reordered/expanded set of column headers (B1:AO1), replacing the
# previous set of headers (B1:AB1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$headers = @(
    "filename",
    "rays_present",
    "approx_split",
    "Diad1_pos",
    "Diad2_pos",
    "HB1_pos",
    "HB2_pos",
    "C13_pos",
    "Diad1_abs_prom",
    "Diad2_abs_prom",
    "HB1_abs_prom",
    "HB2_abs_prom",
    "C13_abs_prom",
    "Mean_abs_HB_prom",
    "Diad2_HB2_abs_prom_ratio",
    "Diad1_HB1_abs_prom_ratio",
    "Diad1_rel_prom",
    "Diad2_rel_prom",
    "HB1_rel_prom",
    "HB2_rel_prom",
    "C13_rel_prom",
    "Diad1_HB1_Valley_prom",
    "Mean_Diad_HB_Valley_prom",
    "Diad1_prom/std_betweendiads",
    "Diad2_prom/std_betweendiads",
    "Av_Diad_prom/std_betweendiads",
    "C13_prom/HB2_prom",
    "Av_Diad_HB_prom_ratio",
    "Left_vs_Right",
    "Diad2_height",
    "HB2_height",
    "C13_height",
    "Diad1_height",
    "HB1_height",
    "Diad1_Median_Bck",
    "Diad2_Median_Bck",
    "C13_HB2_abs_prom_ratio",
    "Diad2_HB2_Valley_prom",
    "HB1_prom/std_betweendiads",
    "HB2_prom/std_betweendiads"
)

$startCol = 2  # column B
$lastCol = $startCol + $headers.Length - 1   # column AO

# Grab the header style already used by the existing header cell (B1) so the
# newly added cells (AC1:AO1) can be formatted the same way as the rest of
# the header row, without disturbing the formatting already present on the
# untouched cells.
$templateCell = $ws.Cells.Item(1, $startCol)

# First, extend formatting for the newly introduced columns (AC:AO) by
# copying the existing header cell's formatting before writing new values.
$oldLastCol = 28  # column AB, the previous last header column
if ($lastCol -gt $oldLastCol) {
    $newRange = $ws.Range($ws.Cells.Item(1, $oldLastCol + 1), $ws.Cells.Item(1, $lastCol))
    $templateCell.Copy() | Out-Null
    $newRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = $false
}

# Now write the full set of header values, in their new order, across
# B1:AO1.
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $headers[$i]
}
